$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.502.99"
$ws.Range("E2").Value = "  -2.81%  "

$ws.Range("D3").Value = "2.470.69"
$ws.Range("E3").Value = "  -2.52%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "311.81"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").Value = "94.55"
$ws.Range("E6").Value = "  -6.17%  "

$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  -3.20%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("D10").Value = "33.62"
$ws.Range("E10").Value = "  -6.71%  "

$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").Value = "  -2.55%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "6.97"
$ws.Range("E13").Value = "  -5.06%  "

$ws.Range("D14").Value = "2.850.00"
$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").Value = "2.480.63"
$ws.Range("E15").Value = "  -3.92%  "

$ws.Range("D16").Value = "14.62"
$ws.Range("E16").Value = "  -8.13%  "

$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  -2.46%  "

$ws.Range("D18").Value = "41.486.81"
$ws.Range("E18").Value = "  -2.78%  "

$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -5.70%  "

$ws.Range("D20").Value = "0.0₃0918"
$ws.Range("E20").Value = "  -3.63%  "

$ws.Range("D21").Value = "11.60"
$ws.Range("E21").Value = "  -5.05%  "

$ws.Range("D22").Value = "69.07"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").Value = "238.11"
$ws.Range("E23").Value = "  -2.27%  "

$ws.Range("D24").Value = "2.78"
$ws.Range("E24").Value = "  -3.92%  "

$ws.Range("E25").Value = "  -4.96%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "24.89"
$ws.Range("E27").Value = "  -4.30%  "

$ws.Range("E28").Value = "  -4.13%  "

$ws.Range("D29").Value = "9.75"
$ws.Range("E29").Value = "  -3.89%  "

$ws.Range("D30").Value = "36.50"
$ws.Range("E30").Value = "  -6.96%  "

$ws.Range("D31").Value = "153.36"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").Value = "5.67"
$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("E34").Value = "  -7.71%  "

$ws.Range("D35").Value = "0.0754"
$ws.Range("E35").Value = "  -4.78%  "

$ws.Range("D36").Value = "3.03"
$ws.Range("E36").Value = "  -4.54%  "

$ws.Range("D37").Value = "17.20"
$ws.Range("E37").Value = "  -6.89%  "

$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  -7.01%  "

$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  -6.37%  "

$ws.Range("E40").Value = "  -3.73%  "

$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -5.08%  "

$ws.Range("D42").Value = "21.41"
$ws.Range("E42").Value = "  -2.52%  "

$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").Value = "1.989.26"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "0.0286"
$ws.Range("E45").Value = "  -4.15%  "

$ws.Range("D46").Value = "3.07"
$ws.Range("E46").Value = "  -7.12%  "

$ws.Range("D47").Value = "8.79"
$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").Value = "77.44"
$ws.Range("E48").Value = "  -4.71%  "

$ws.Range("D49").Value = "97.93"
$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("D50").Value = "69.39"
$ws.Range("E50").Value = "  -4.34%  "

$ws.Range("D51").Value = "0.181"
$ws.Range("E51").Value = "  -5.89%  "
